$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").Value = "8OJEKV"
$ws.Range("B60").Value = "Film de fusor Brother"
$ws.Range("C60").Value = "HL 5440 5445 5450 5452 5455 5470 5472 5580 5585 5590 5595 6180 6182 6200 6300 7810, DCP 8110 8112 8150 8152 8155 8510, MFC 8510 8512 8515 8520 8710 8712 8910 8912 8950 8952"
$ws.Range("D60").Value = 100000
$ws.Range("E60").Value = 200000
$ws.Range("F60").Value = 3
$ws.Range("G60").Value = 1
$ws.Range("H60").Formula = "=(E60-D60)*G60"
$ws.Range("I60").Formula = "=D60*F60"
$ws.Range("J60").Value = 300000
